$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.976.55"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "3.170.50"
$ws.Range("E4").Value = "  -0.01%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'579.73"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  +4.38%  "
$__style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'151.00"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  +6.81%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.170.04"
$ws.Range("E8").Value = "  +4.17%  "
$__style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.530"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("E10").Value = "  +6.51%  "
$__style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'6.20"
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = "  +0.43%  "
$__style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.501"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = "  +3.61%  "
$__style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = "  +18.89%  "
$__style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'37.56"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = "  +5.98%  "
$ws.Range("D15").Value = "3.687.49"
$ws.Range("E15").Value = "  +4.21%  "
$ws.Range("D16").Value = "65.051.71"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$__style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'7.17"
$ws.Range("D17").Style = $__style
$ws.Range("E17").Value = "  +6.24%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.160.46"
$ws.Range("E18").Value = "  +4.13%  "
$ws.Range("E19").Value = "  +1.44%  "
$__style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'511.29"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = "  +6.93%  "
$__style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'14.87"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = "  +5.96%  "
$__style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'0.727"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = "  +6.85%  "
$__style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'15.42"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  +6.59%  "
$__style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'7.81"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = "  +3.89%  "
$__style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'85.07"
$ws.Range("D25").Style = $__style
$ws.Range("E26").Value = "  +0.14%  "
$__style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'9.09"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = "  +12.77%  "
$__style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'2.93"
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = "  +5.23%  "
$__style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'2.19"
$ws.Range("D29").Style = $__style
$ws.Range("E29").Value = "  +8.30%  "
$ws.Range("B30").Value = "Stacks"
$ws.Range("C30").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$__style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'2.81"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = "  +15.72%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$__style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'27.81"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  +6.71%  "
$__style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  +4.03%  "
$__style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'6.35"
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = "  +12.22%  "
$__style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'6.59"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  +6.50%  "
$__style = $ws.Range("D36").Style
$ws.Range("D36").Value = "'55.67"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  +1.36%  "
$__style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.0901"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = "  +11.17%  "
$__style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'474.13"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  +8.19%  "
$ws.Range("E39").Value = "  +14.25%  "
$__style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.0421"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = "  +3.38%  "
$__style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'8.65"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = "  +4.67%  "
$ws.Range("D42").Value = "3.066.02"
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("E43").Value = "  +1.62%  "
$__style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'0.284"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +6.01%  "
$__style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'2.44"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = "  +9.58%  "
$__style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'29.39"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = "  +6.03%  "
$ws.Range("D47").Value = "0.0₃0609"
$ws.Range("E47").Value = "  +19.51%  "
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("E50").Value = "  +8.57%  "
$__style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'120.26"
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = "  +2.09%  "
